$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on Price (D) cells whose new value would otherwise
# be auto-recognized by Excel as a number, so they remain text like the source data.
$textPriceCells = @(
    "D5",
    "D7",
    "D9",
    "D11",
    "D12",
    "D14",
    "D15",
    "D19",
    "D20",
    "D21",
    "D22",
    "D27",
    "D29",
    "D30",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D46",
    "D48",
    "D50"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
$ws.Range("D2").Value = "43.882.38"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.348.01"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "239.88"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("D7").Value = "73.66"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "60.97"
$ws.Range("E11").Value = "  +6.58%  "
$ws.Range("D12").Value = "33.36"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "7.26"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "16.17"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "2.351.87"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "43.770.06"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "6.62"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "77.88"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "252.64"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  -3.30%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "175.69"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("D30").Value = "22.21"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").Value = "0.0743"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").Value = "5.36"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "5.05"
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("D36").Value = "3.78"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").Value = "6.43"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("D40").Value = "5.39"
$ws.Range("E40").Value = "  +12.32%  "
$ws.Range("D41").Value = "65.72"
$ws.Range("E41").Value = "  +15.71%  "
$ws.Range("D42").Value = "19.35"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "9.16"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("E44").Value = "  -3.04%  "
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("D50").Value = "98.10"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("E51").Value = "  +2.79%  "
